$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BON DE COMMANDE 1")
Write-Host "=== BEFORE ==="
Write-Host "B19 Font.Color" $ws.Range("B19").Font.Color
Write-Host "B19 Font.Underline" $ws.Range("B19").Font.Underline
Write-Host "B19 Interior.Color" $ws.Range("B19").Interior.Color
Write-Host "B19 Borders(7).LineStyle" $ws.Range("B19").Borders.Item(7).LineStyle
Write-Host "B19 Borders(7).Weight" $ws.Range("B19").Borders.Item(7).Weight
Write-Host "B19 Borders(10).LineStyle" $ws.Range("B19").Borders.Item(10).LineStyle
